# Fill in the "Status" grid (columns F:I) for the task rows with "Done",
# mirroring a manual pass where the checklist cells got marked complete.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "Done"
$ws.Range("I2").Value = "Done"
$ws.Range("F3:I3").Value = "Done"
$ws.Range("F4:I4").Value = "Done"
$ws.Range("F6").Value = "Done"

# Leave the new selection where the user's cursor ended up.
$ws.Range("F3").Select()
